$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Update the quantity ordered for the second row (Ipoh Coffee) to an amount
# that exceeds what is available, and record the resulting failure status
# and note - mirroring the "check for available quantity" activity change.
$ws.Range("B2").Value = 2900
$ws.Range("C2").Value = "Failed"
$ws.Range("D2").Value = "Quantity '2900' was unavailable"

# Let the row grow to fit the wrapped note text (same behavior already
# visible on rows 9/10 for the other failure notes).
$ws.Rows.Item(2).RowHeight = 29

# Move the active selection to the cell that was just edited.
$ws.Range("B2").Select()
